$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.777.22"
$ws.Range("D3").Value = "3.810.46"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.02"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.92"
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("D7").Value = "3.808.19"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("E10").Value = "  -4.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.28"
$ws.Range("E11").Value = "  -6.19%  "
$ws.Range("E12").Value = "  -3.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.74"
$ws.Range("E13").Value = "  -3.58%  "
$ws.Range("E14").Value = "  -3.50%  "
$ws.Range("D15").Value = "4.444.39"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "3.808.54"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "67.784.58"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.26"
$ws.Range("E18").Value = "  -3.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.115"
$ws.Range("E19").Value = "  -4.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.95"
$ws.Range("E20").Value = "  +3.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "493.47"
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.25"
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.13"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("E25").Value = "  +7.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.38"
$ws.Range("E26").Value = "  -5.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.33"
$ws.Range("E27").Value = "  -3.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.27"
$ws.Range("E28").Value = "  -3.70%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.99"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.22"
$ws.Range("E31").Value = "  +8.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.44"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.79"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("E34").Value = "  -4.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.83"
$ws.Range("E37").Value = "  -5.20%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "463.58"
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.133"
$ws.Range("E39").Value = "  -4.13%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.330"
$ws.Range("E40").Value = "  -3.50%  "
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("E42").Value = "  -2.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.86"
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.45"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.24"
$ws.Range("E45").Value = "  -7.43%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.852.93"
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "138.99"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.96"
$ws.Range("E50").Value = "  -4.71%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.75"
$ws.Range("E51").Value = "  +6.97%  "
